$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis), (samp:sample)-->(c)
WHERE s.clinical_study_designation IN ['OSA01'] and s.clinical_study_type in ['Genomics']  and demo.sex in ['Female'] and samp.sample_site in ['Bone'] and samp.specific_sample_pathology in ['Osteoblastic Osteosarcoma']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  
       coalesce(c.case_id, '') AS `Case ID`,
       coalesce(s.clinical_study_designation, '') AS `Study Code`,
       coalesce(s.clinical_study_type, '') AS  `Study Type`,
       coalesce(demo.breed, '') AS Breed ,
       coalesce(diag.disease_term, '') AS Diagnosis ,
       coalesce(diag.stage_of_disease, '') AS `Stage of Disease`,
       CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
 coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`

'@
$samplesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis), (samp:sample)-->(c)
WHERE s.clinical_study_designation  IN ['OSA01'] and s.clinical_study_type in ['Genomics']  and demo.sex in ['Female'] and samp.sample_site in ['Bone'] and samp.specific_sample_pathology in ['Osteoblastic Osteosarcoma']
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
'@
$studyFilesQuery = @'
MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
MATCH (f)-[*]->(samp:sample)
WHERE s.clinical_study_designation  IN ['OSA01'] and s.clinical_study_type in ['Genomics']  and demo.sex in ['Female'] and samp.sample_site in ['Bone'] and samp.specific_sample_pathology in ['Osteoblastic Osteosarcoma']
WITH DISTINCT f, s, c, demo, diag
WITH
        f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH    
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH    
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
'@
$statQuery = @'
MATCH (s:study)
MATCH (demo:demographic) 
MATCH (diag:diagnosis)
MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis), (samp:sample)-->(c)
WHERE s.clinical_study_designation  IN ['OSA01'] and s.clinical_study_type in ['Genomics']  and demo.sex in ['Female'] and samp.sample_site in ['Bone'] and samp.specific_sample_pathology in ['Osteoblastic Osteosarcoma']
    
OPTIONAL MATCH (s)<-[:member_of]-(c:case)
OPTIONAL MATCH (c)<-[:of_case]-(samp:sample)<-[:of_sample]-(f:file)
OPTIONAL MATCH (sf:file)-->(s)

RETURN 
count(DISTINCT(p)) as Programs,
count(DISTINCT(s.clinical_study_designation)) as Studies,
count(DISTINCT(c.case_id)) as Cases,
count(DISTINCT(samp)) as Samples , 
count(DISTINCT(f)) as CaseFiles , count(distinct sf) AS `Study Files`
'@
$filesQuery = @'
MATCH (c:case)<--(demo:demographic)
MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c)
MATCH (f:file)-[*]->(c)
MATCH (f)-->(parent)
WHERE s.clinical_study_designation  IN ['OSA01'] and s.clinical_study_type in ['Genomics']  and demo.sex in ['Female']
 and samp.sample_site in ['Bone'] and samp.specific_sample_pathology in ['Osteoblastic Osteosarcoma']
OPTIONAL MATCH (f)-[*]->(smpl:sample)
WITH DISTINCT f, smpl, parent, c, diag, demo, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, smpl, parent, c, diag, demo, s,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, smpl, parent, c, diag, demo, s, unit,
        round(factor * value)/factor AS size
RETURN distinct
        coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_type, '') AS `File Type`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(smpl.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
'@

$ws.Range("B2").Value = $casesQuery
$ws.Range("B3").Value = $samplesQuery
$ws.Range("B5").Value = $studyFilesQuery
$ws.Range("C2:C5").Value = $statQuery
$ws.Range("B4").Value = $filesQuery

$ws.Rows(4).RowHeight = 409.5
$ws.Rows(5).RowHeight = 409.5

$ws.Range("D5").Select() | Out-Null

Write-Host "Edit applied successfully"
